$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header and data for column G:
# Row 1 header keeps text "MSRP_VERSION" (shared string index just gets reindexed automatically by the writer)
$ws.Range("G1").Value = "MSRP_VERSION"
# Row 2: MSRP_2000 -> MSRP_2000_SELECT
$ws.Range("G2").Value = "MSRP_2000_SELECT"
# Row 3: MSRP_2000_CA_SELECT -> MSRP_2000_SELECT
$ws.Range("G3").Value = "MSRP_2000_SELECT"

# Update the selected/active cell shown when the workbook is opened
$ws.Range("J12").Select()
